$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Ace"
$ws.Cells.Item(2,3).Value2 = "Bdkrb2"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 53.74035766666666
$ws.Cells.Item(2,8).Value2 = 161.221073
$ws.Cells.Item(2,9).Value2 = 0.4609242213031295
$ws.Cells.Item(2,10).Value2 = 0.4621087068901442
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 0.9242423333333334
$ws.Cells.Item(2,14).Value2 = 2.772727
$ws.Cells.Item(2,15).Value2 = 0.7379131210038523
$ws.Cells.Item(2,16).Value2 = 0.8085498937233963
$ws.Cells.Item(2,17).Value2 = 49.66911356400789
$ws.Cells.Item(2,18).Value2 = 447.022022076071
$ws.Cells.Item(2,19).Value2 = 0.3401220306880626
$ws.Cells.Item(2,20).Value2 = 0.3736379458446822
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Ace"
$ws.Cells.Item(3,3).Value2 = "Bdkrb2"
$ws.Cells.Item(3,4).Value2 = "MuSCs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 53.74035766666666
$ws.Cells.Item(3,8).Value2 = 161.221073
$ws.Cells.Item(3,9).Value2 = 0.4609242213031295
$ws.Cells.Item(3,10).Value2 = 0.4621087068901442
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 0.328266
$ws.Cells.Item(3,14).Value2 = 0.656532
$ws.Cells.Item(3,15).Value2 = 0.2620868789961477
$ws.Cells.Item(3,16).Value2 = 0.1914501062766038
$ws.Cells.Item(3,17).Value2 = 17.641132249806
$ws.Cells.Item(3,18).Value2 = 105.846793498836
$ws.Cells.Item(3,19).Value2 = 0.1208021906150669
$ws.Cells.Item(3,20).Value2 = 0.08847076104546206
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Ace"
$ws.Cells.Item(4,3).Value2 = "Bdkrb2"
$ws.Cells.Item(4,4).Value2 = "ECs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 54.64926166666667
$ws.Cells.Item(4,8).Value2 = 163.947785
$ws.Cells.Item(4,9).Value2 = 0.4687197754570081
$ws.Cells.Item(4,10).Value2 = 0.4699242941017605
$ws.Cells.Item(4,11).Value2 = 2
$ws.Cells.Item(4,12).Value2 = 0.6666666666666666
$ws.Cells.Item(4,13).Value2 = 0.9242423333333334
$ws.Cells.Item(4,14).Value2 = 2.772727
$ws.Cells.Item(4,15).Value2 = 0.7379131210038523
$ws.Cells.Item(4,16).Value2 = 0.8085498937233963
$ws.Cells.Item(4,17).Value2 = 50.50916111774389
$ws.Cells.Item(4,18).Value2 = 454.5824500596951
$ws.Cells.Item(4,19).Value2 = 0.3458744723837057
$ws.Cells.Item(4,20).Value2 = 0.3799572380540204
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Ace"
$ws.Cells.Item(5,3).Value2 = "Bdkrb2"
$ws.Cells.Item(5,4).Value2 = "MuSCs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 54.64926166666667
$ws.Cells.Item(5,8).Value2 = 163.947785
$ws.Cells.Item(5,9).Value2 = 0.4687197754570081
$ws.Cells.Item(5,10).Value2 = 0.4699242941017605
$ws.Cells.Item(5,11).Value2 = 2
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 0.328266
$ws.Cells.Item(5,14).Value2 = 0.656532
$ws.Cells.Item(5,15).Value2 = 0.2620868789961477
$ws.Cells.Item(5,16).Value2 = 0.1914501062766038
$ws.Cells.Item(5,17).Value2 = 17.93949453027
$ws.Cells.Item(5,18).Value2 = 107.63696718162
$ws.Cells.Item(5,19).Value2 = 0.1228453030733024
$ws.Cells.Item(5,20).Value2 = 0.08996705604774005
$ws.Cells.Item(6,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value2 = "Ace"
$ws.Cells.Item(6,3).Value2 = "Bdkrb2"
$ws.Cells.Item(6,4).Value2 = "ECs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 5.296475666666667
$ws.Cells.Item(6,8).Value2 = 15.889427
$ws.Cells.Item(6,9).Value2 = 0.04542719900473509
$ws.Cells.Item(6,10).Value2 = 0.04554393806940699
$ws.Cells.Item(6,11).Value2 = 2
$ws.Cells.Item(6,12).Value2 = 0.6666666666666666
$ws.Cells.Item(6,13).Value2 = 0.9242423333333334
$ws.Cells.Item(6,14).Value2 = 2.772727
$ws.Cells.Item(6,15).Value2 = 0.7379131210038523
$ws.Cells.Item(6,16).Value2 = 0.8085498937233963
$ws.Cells.Item(6,17).Value2 = 4.895227028603222
$ws.Cells.Item(6,18).Value2 = 44.057043257429
$ws.Cells.Item(6,19).Value2 = 0.03352132619604716
$ws.Cells.Item(6,20).Value2 = 0.03682454628576397
$ws.Cells.Item(7,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value2 = "Ace"
$ws.Cells.Item(7,3).Value2 = "Bdkrb2"
$ws.Cells.Item(7,4).Value2 = "MuSCs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 5.296475666666667
$ws.Cells.Item(7,8).Value2 = 15.889427
$ws.Cells.Item(7,9).Value2 = 0.04542719900473509
$ws.Cells.Item(7,10).Value2 = 0.04554393806940699
$ws.Cells.Item(7,11).Value2 = 2
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 0.328266
$ws.Cells.Item(7,14).Value2 = 0.656532
$ws.Cells.Item(7,15).Value2 = 0.2620868789961477
$ws.Cells.Item(7,16).Value2 = 0.1914501062766038
$ws.Cells.Item(7,17).Value2 = 1.738652881194
$ws.Cells.Item(7,18).Value2 = 10.431917287164
$ws.Cells.Item(7,19).Value2 = 0.01190587280868793
$ws.Cells.Item(7,20).Value2 = 0.008719391783643029
$ws.Cells.Item(8,1).Value2 = "MuSCs"
$ws.Cells.Item(8,2).Value2 = "Ace"
$ws.Cells.Item(8,3).Value2 = "Bdkrb2"
$ws.Cells.Item(8,4).Value2 = "ECs"
$ws.Cells.Item(8,5).Value2 = 2
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 0.896557
$ws.Cells.Item(8,8).Value2 = 1.793114
$ws.Cells.Item(8,9).Value2 = 0.007689655503264204
$ws.Cells.Item(8,10).Value2 = 0.005139610948046563
$ws.Cells.Item(8,11).Value2 = 2
$ws.Cells.Item(8,12).Value2 = 0.6666666666666666
$ws.Cells.Item(8,13).Value2 = 0.9242423333333334
$ws.Cells.Item(8,14).Value2 = 2.772727
$ws.Cells.Item(8,15).Value2 = 0.7379131210038523
$ws.Cells.Item(8,16).Value2 = 0.8085498937233963
$ws.Cells.Item(8,17).Value2 = 0.8286359336463335
$ws.Cells.Item(8,18).Value2 = 4.971815601878
$ws.Cells.Item(8,19).Value2 = 0.005674297691858138
$ws.Cells.Item(8,20).Value2 = 0.004155631885822653
$ws.Cells.Item(9,1).Value2 = "MuSCs"
$ws.Cells.Item(9,2).Value2 = "Ace"
$ws.Cells.Item(9,3).Value2 = "Bdkrb2"
$ws.Cells.Item(9,4).Value2 = "MuSCs"
$ws.Cells.Item(9,5).Value2 = 2
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 0.896557
$ws.Cells.Item(9,8).Value2 = 1.793114
$ws.Cells.Item(9,9).Value2 = 0.007689655503264204
$ws.Cells.Item(9,10).Value2 = 0.005139610948046563
$ws.Cells.Item(9,11).Value2 = 2
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 0.328266
$ws.Cells.Item(9,14).Value2 = 0.656532
$ws.Cells.Item(9,15).Value2 = 0.2620868789961477
$ws.Cells.Item(9,16).Value2 = 0.1914501062766038
$ws.Cells.Item(9,17).Value2 = 0.294309180162
$ws.Cells.Item(9,18).Value2 = 1.177236720648
$ws.Cells.Item(9,19).Value2 = 0.002015357811406067
$ws.Cells.Item(9,20).Value2 = 0.0009839790622239108
$ws.Cells.Item(10,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(10,2).Value2 = "Ace"
$ws.Cells.Item(10,3).Value2 = "Bdkrb2"
$ws.Cells.Item(10,4).Value2 = "ECs"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 2.009957333333333
$ws.Cells.Item(10,8).Value2 = 6.029872
$ws.Cells.Item(10,9).Value2 = 0.01723914873186302
$ws.Cells.Item(10,10).Value2 = 0.01728344999064166
$ws.Cells.Item(10,11).Value2 = 2
$ws.Cells.Item(10,12).Value2 = 0.6666666666666666
$ws.Cells.Item(10,13).Value2 = 0.9242423333333334
$ws.Cells.Item(10,14).Value2 = 2.772727
$ws.Cells.Item(10,15).Value2 = 0.7379131210038523
$ws.Cells.Item(10,16).Value2 = 0.8085498937233963
$ws.Cells.Item(10,17).Value2 = 1.857687655660445
$ws.Cells.Item(10,18).Value2 = 16.719188900944
$ws.Cells.Item(10,19).Value2 = 0.01272099404417864
$ws.Cells.Item(10,20).Value2 = 0.01397453165310695
$ws.Cells.Item(11,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(11,2).Value2 = "Ace"
$ws.Cells.Item(11,3).Value2 = "Bdkrb2"
$ws.Cells.Item(11,4).Value2 = "MuSCs"
$ws.Cells.Item(11,5).Value2 = 3
$ws.Cells.Item(11,6).Value2 = 1
$ws.Cells.Item(11,7).Value2 = 2.009957333333333
$ws.Cells.Item(11,8).Value2 = 6.029872
$ws.Cells.Item(11,9).Value2 = 0.01723914873186302
$ws.Cells.Item(11,10).Value2 = 0.01728344999064166
$ws.Cells.Item(11,11).Value2 = 2
$ws.Cells.Item(11,12).Value2 = 1
$ws.Cells.Item(11,13).Value2 = 0.328266
$ws.Cells.Item(11,14).Value2 = 0.656532
$ws.Cells.Item(11,15).Value2 = 0.2620868789961477
$ws.Cells.Item(11,16).Value2 = 0.1914501062766038
$ws.Cells.Item(11,17).Value2 = 0.659800653984
$ws.Cells.Item(11,18).Value2 = 3.958803923904
$ws.Cells.Item(11,19).Value2 = 0.004518154687684376
$ws.Cells.Item(11,20).Value2 = 0.003308918337534712

